# The deck ships two DrawingML theme parts:
#   ppt/theme/theme1.xml -> linked only from the Notes Master ("Office Theme" colours)
#   ppt/theme/theme2.xml -> linked from the (one and only) Slide Master / Design
#                            ("Integral" colours) -- this is the theme that is actually
#                            visible throughout the presentation.
#
# The authored edit swaps the two parts' contents (Integral <-> Office Theme). The
# font scheme and format (fill/line/effect) scheme are identical between the two
# themes, so the only real difference is the 12-slot colour scheme. We reproduce
# that swap on the live/visible design via the ThemeColorScheme object, which is
# the supported COM surface for rewriting a design's colour scheme.

$p = $ppt.ActivePresentation

# Target palette = the colours the "Office Theme" part (theme1.xml) currently has;
# after the swap this is what the active design (theme2.xml) must become.
$officeThemeColors = @(
    "000000",  # 1  dk1
    "FFFFFF",  # 2  lt1
    "44546A",  # 3  dk2
    "E7E6E6",  # 4  lt2
    "5B9BD5",  # 5  accent1
    "ED7D31",  # 6  accent2
    "A5A5A5",  # 7  accent3
    "FFC000",  # 8  accent4
    "4472C4",  # 9  accent5
    "70AD47",  # 10 accent6
    "0563C1",  # 11 hlink
    "954F72"   # 12 folHlink
)

$slide = $p.Slides.Item(1)
$themeColors = $slide.ThemeColorScheme

for ($i = 1; $i -le $themeColors.Count; $i++) {
    $hex = $officeThemeColors[$i - 1]
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $themeColors.Colors($i).RGB = $r + ($g * 256) + ($b * 65536)
}
